$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("wind")

$ws.Range("P4").Value = 3
$ws.Range("P5").Value = 2

$ws.Range("P13").Value = 4
$ws.Range("P14").Value = 5

$ws.Range("P15").Value = 2
$ws.Range("P16").Value = 3
$ws.Range("P17").Value = 1

$ws.Range("P19").Value = 2
$ws.Range("P20").Value = 1

$ws.Range("P27").Value = 1
$ws.Range("P28").Value = 2

$ws.Range("P47").Value = 1
$ws.Range("P48").Value = 2
